# Generate Report for Archive
# The three localized files that were previously "Ready for handoff" have now
# progressed to "In Translation". Update the Overview sheet (per-language
# status columns) and each language sheet's Status column accordingly.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E5").Value = "In Translation"
$overview.Range("F5").Value = "In Translation"
$overview.Range("E6").Value = "In Translation"
$overview.Range("F6").Value = "In Translation"
$overview.Range("E7").Value = "In Translation"
$overview.Range("F7").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C5").Value = "In Translation"
$zhcn.Range("C6").Value = "In Translation"
$zhcn.Range("C7").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C5").Value = "In Translation"
$dede.Range("C6").Value = "In Translation"
$dede.Range("C7").Value = "In Translation"
